# Apply data-cleaning fixes to the Arizona 2019 MCAS sheet:
#  1. Rename header columns to the cleaned technical names.
#  2. Title-case the state / municipality name text (e.g. "de" -> "De", "del" -> "Del").
#  3. Remove the trailing metadata / footnote rows that were appended after the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the footer/metadata rows (rows 1241-1245) -------------------
# These hold free-text notes below the actual data table and are dropped entirely.
$ws.Range("A1241:A1245").EntireRow.Delete()

# --- 2. Title-case the textual columns (A = state, B = municipality) -------
# Excel's UsedRange now ends at row 1239 after the footer rows were removed.
$lastRow = $ws.UsedRange.Rows.Count

$rng = $ws.Range("A2:B" + $lastRow)
$vals = $rng.Value2

for ($i = 1; $i -le $vals.GetLength(0); $i++) {
    for ($j = 1; $j -le $vals.GetLength(1); $j++) {
        $v = $vals[$i, $j]
        if ($v -ne $null -and $v.GetType().Name -eq "String") {
            $words = $v.Split(" ")
            $fixedWords = @()
            foreach ($w in $words) {
                if ($w.Length -gt 0) {
                    $fixedWords += [string]::Concat($w.Substring(0, 1).ToUpper(), $w.Substring(1))
                } else {
                    $fixedWords += $w
                }
            }
            $vals[$i, $j] = [string]::Join(" ", $fixedWords)
        }
    }
}

$rng.Value2 = $vals

# --- 3. Rename the header row to the cleaned column names -------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"
